$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J2:J11
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: summary stats
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style for B14:B17 - bold size 12 font, vertical center alignment
$c = $ws.Range("B14")
$c.Font.Bold = $true
$c.Font.Size = 12
$c.VerticalAlignment = -4108  # xlCenter
$c.Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row height for rows 14-17 (15.6)
$ws.Range("A14:A17").RowHeight = 15.6

# Page setup: paper size 9 (A4), portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("J2:J12").Select()
